$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the existing row 224, shifting the old
# rows 224-267 down to 227-270.
$ws.Range("A224:A226").EntireRow.Insert()

# --- New row 224 ---
$ws.Cells.Item(224, 1).Value = 9
$ws.Cells.Item(224, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(224, 3).Value = "Metropolitana"
$ws.Cells.Item(224, 4).Value = 44644
$ws.Cells.Item(224, 5).Value = 13
$ws.Cells.Item(224, 6).Value = 100112021
$ws.Cells.Item(224, 7).Value = "Ají"
$ws.Cells.Item(224, 8).Value = "Inferno"
$ws.Cells.Item(224, 9).Value = "Primera"
$ws.Cells.Item(224, 10).Value = 34
$ws.Cells.Item(224, 11).Value = 15000
$ws.Cells.Item(224, 12).Value = 15000
$ws.Cells.Item(224, 13).Value = 15000
$ws.Cells.Item(224, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(224, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(224, 16).Value = 1250
$ws.Cells.Item(224, 17).Value = 12
$ws.Cells.Item(224, 18).Value = "Hortaliza"

# --- New row 225 ---
$ws.Cells.Item(225, 1).Value = 9
$ws.Cells.Item(225, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(225, 3).Value = "Metropolitana"
$ws.Cells.Item(225, 4).Value = 44644
$ws.Cells.Item(225, 5).Value = 13
$ws.Cells.Item(225, 6).Value = 100112021
$ws.Cells.Item(225, 7).Value = "Ají"
$ws.Cells.Item(225, 8).Value = "Inferno"
$ws.Cells.Item(225, 9).Value = "Primera"
$ws.Cells.Item(225, 10).Value = 45
$ws.Cells.Item(225, 11).Value = 20000
$ws.Cells.Item(225, 12).Value = 20000
$ws.Cells.Item(225, 13).Value = 20000
$ws.Cells.Item(225, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(225, 15).Value = "Limache"
$ws.Cells.Item(225, 16).Value = 1333
$ws.Cells.Item(225, 17).Value = 15
$ws.Cells.Item(225, 18).Value = "Hortaliza"

# --- New row 226 ---
$ws.Cells.Item(226, 1).Value = 9
$ws.Cells.Item(226, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(226, 3).Value = "Metropolitana"
$ws.Cells.Item(226, 4).Value = 44644
$ws.Cells.Item(226, 5).Value = 13
$ws.Cells.Item(226, 6).Value = 100112021
$ws.Cells.Item(226, 7).Value = "Ají"
$ws.Cells.Item(226, 8).Value = "Inferno"
$ws.Cells.Item(226, 9).Value = "Segunda"
$ws.Cells.Item(226, 10).Value = 30
$ws.Cells.Item(226, 11).Value = 18000
$ws.Cells.Item(226, 12).Value = 18000
$ws.Cells.Item(226, 13).Value = 18000
$ws.Cells.Item(226, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(226, 15).Value = "Limache"
$ws.Cells.Item(226, 16).Value = 1200
$ws.Cells.Item(226, 17).Value = 15
$ws.Cells.Item(226, 18).Value = "Hortaliza"
